# Weekly update: a new Coliflor price record (week of 2023-01-05) was
# reported for "Vega Monumental Concepción". It slots in chronologically
# just after row 342 (2022-08-?? / serial 44813) and before the former
# row 343 (2021-12-15 / serial 44545), so insert a new row at 343 and
# shift the existing rows 343:363 down to 344:364.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(343).Insert()

$ws.Range("A343").Value = 11
$ws.Range("B343").Value = "Vega Monumental Concepción"
$ws.Range("C343").Value = "Bíobío"
$ws.Range("D343").Value = 44931
$ws.Range("E343").Value = 8
$ws.Range("F343").Value = 100112008
$ws.Range("G343").Value = "Coliflor"
$ws.Range("H343").Value = "Sin especificar"
$ws.Range("I343").Value = "Primera"
$ws.Range("J343").Value = 2000
$ws.Range("K343").Value = 800
$ws.Range("L343").Value = 850
$ws.Range("M343").Value = 825
$ws.Range("N343").Value = "`$/unidad"
$ws.Range("O343").Value = "Región Metropolitana"
$ws.Range("P343").Value = 825
$ws.Range("Q343").Value = 1
$ws.Range("R343").Value = "Hortaliza"
